$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

${ws}.Range('D2').Value = '45.795.88'
${ws}.Range('E2').Value = '  +6.50%  '
${ws}.Range('D3').Value = '2.415.79'
${ws}.Range('E3').Value = '  +5.09%  '
${ws}.Range('E4').Value = '  +0.01%  '
${ws}.Range('D5').NumberFormat = '@'
${ws}.Range('D5').Value = '116.43'
${ws}.Range('E5').Value = '  +11.60%  '
${ws}.Range('D6').NumberFormat = '@'
${ws}.Range('D6').Value = '319.27'
${ws}.Range('E6').Value = '  +2.45%  '
${ws}.Range('D7').NumberFormat = '@'
${ws}.Range('D7').Value = '0.636'
${ws}.Range('E7').Value = '  +1.98%  '
${ws}.Range('E8').Value = '  -0.09%  '
${ws}.Range('D9').NumberFormat = '@'
${ws}.Range('D9').Value = '0.631'
${ws}.Range('E9').Value = '  +4.43%  '
${ws}.Range('D10').NumberFormat = '@'
${ws}.Range('D10').Value = '43.36'
${ws}.Range('E10').Value = '  +9.45%  '
${ws}.Range('E11').Value = '  +4.39%  '
${ws}.Range('D12').NumberFormat = '@'
${ws}.Range('D12').Value = '8.77'
${ws}.Range('E12').Value = '  +6.06%  '
${ws}.Range('E13').Value = '  +4.02%  '
${ws}.Range('E14').Value = '  +2.17%  '
${ws}.Range('D15').NumberFormat = '@'
${ws}.Range('D15').Value = '15.94'
${ws}.Range('E15').Value = '  +3.98%  '
${ws}.Range('D16').Value = '2.786.53'
${ws}.Range('E16').Value = '  +5.36%  '
${ws}.Range('D17').Value = '2.420.67'
${ws}.Range('E17').Value = '  +5.38%  '
${ws}.Range('D18').Value = '45.820.66'
${ws}.Range('E18').Value = '  +7.04%  '
${ws}.Range('D19').NumberFormat = '@'
${ws}.Range('D19').Value = '7.62'
${ws}.Range('E19').Value = '  +4.20%  '
${ws}.Range('E20').Value = '  +4.22%  '
${ws}.Range('D21').NumberFormat = '@'
${ws}.Range('D21').Value = '13.46'
${ws}.Range('E21').Value = '  +0.16%  '
${ws}.Range('D22').NumberFormat = '@'
${ws}.Range('D22').Value = '75.18'
${ws}.Range('E22').Value = '  +2.41%  '
${ws}.Range('E23').Value = '  +4.34%  '
${ws}.Range('D24').NumberFormat = '@'
${ws}.Range('D24').Value = '268.19'
${ws}.Range('E24').Value = '  +0.20%  '
${ws}.Range('D25').NumberFormat = '@'
${ws}.Range('D25').Value = '2.40'
${ws}.Range('E26').Value = '  -0.58%  '
${ws}.Range('D27').NumberFormat = '@'
${ws}.Range('D27').Value = '7.69'
${ws}.Range('E27').Value = '  +8.09%  '
${ws}.Range('D28').NumberFormat = '@'
${ws}.Range('D28').Value = '11.41'
${ws}.Range('E28').Value = '  +5.52%  '
${ws}.Range('E29').Value = '  +2.62%  '
${ws}.Range('D30').NumberFormat = '@'
${ws}.Range('D30').Value = '40.21'
${ws}.Range('E30').Value = '  +11.21%  '
${ws}.Range('D31').NumberFormat = '@'
${ws}.Range('D31').Value = '23.04'
${ws}.Range('E31').Value = '  +3.06%  '
${ws}.Range('D32').NumberFormat = '@'
${ws}.Range('D32').Value = '0.0975'
${ws}.Range('E32').Value = '  +14.30%  '
${ws}.Range('D33').NumberFormat = '@'
${ws}.Range('D33').Value = '173.91'
${ws}.Range('E33').Value = '  +5.42%  '
${ws}.Range('E34').Value = '  +14.33%  '
${ws}.Range('B35').Value = 'RenderToken'
${ws}.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
${ws}.Range('D35').NumberFormat = '@'
${ws}.Range('D35').Value = '5.02'
${ws}.Range('E35').Value = '  +10.24%  '
${ws}.Range('B36').Value = 'Stellar'
${ws}.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
${ws}.Range('D36').NumberFormat = '@'
${ws}.Range('D36').Value = '0.132'
${ws}.Range('E36').Value = '  +2.14%  '
${ws}.Range('E37').Value = '  +7.28%  '
${ws}.Range('D38').NumberFormat = '@'
${ws}.Range('D38').Value = '4.26'
${ws}.Range('E38').Value = '  +17.26%  '
${ws}.Range('D39').NumberFormat = '@'
${ws}.Range('D39').Value = '3.17'
${ws}.Range('E39').Value = '  +12.28%  '
${ws}.Range('D40').NumberFormat = '@'
${ws}.Range('D40').Value = '0.0366'
${ws}.Range('E40').Value = '  +5.63%  '
${ws}.Range('D41').NumberFormat = '@'
${ws}.Range('D41').Value = '1.81'
${ws}.Range('E41').Value = '  +14.78%  '
${ws}.Range('D42').NumberFormat = '@'
${ws}.Range('D42').Value = '102.34'
${ws}.Range('E42').Value = '  -5.57%  '
${ws}.Range('D43').NumberFormat = '@'
${ws}.Range('D43').Value = '13.66'
${ws}.Range('E43').Value = '  +12.55%  '
${ws}.Range('E44').Value = '  +5.83%  '
${ws}.Range('D45').NumberFormat = '@'
${ws}.Range('D45').Value = '72.56'
${ws}.Range('E45').Value = '  +2.05%  '
${ws}.Range('E46').Value = '  -0.49%  '
${ws}.Range('D47').NumberFormat = '@'
${ws}.Range('D47').Value = '5.86'
${ws}.Range('E47').Value = '  +13.91%  '
${ws}.Range('D48').NumberFormat = '@'
${ws}.Range('D48').Value = '117.55'
${ws}.Range('E48').Value = '  +6.32%  '
${ws}.Range('E49').Value = '  +18.28%  '
${ws}.Range('B50').Value = 'ordi'
${ws}.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
${ws}.Range('D50').NumberFormat = '@'
${ws}.Range('D50').Value = '81.26'
${ws}.Range('E50').Value = '  +4.62%  '
${ws}.Range('B51').Value = 'FraxShare'
${ws}.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
${ws}.Range('D51').NumberFormat = '@'
${ws}.Range('D51').Value = '9.52'
${ws}.Range('E51').Value = '  +9.96%  '
